$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.891.75"
$ws.Range("E2").Value = "  +0.18%  "

$ws.Range("D3").Value = "2.635.80"
$ws.Range("E3").Value = "  -0.06%  "

$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").Value = "'596.02"
$ws.Range("E5").Value = "  -1.46%  "

$ws.Range("D6").Value = "'152.31"
$ws.Range("E6").Value = "  -1.96%  "

$ws.Range("E7").Value = "  +0.05%  "

$ws.Range("D8").Value = "'0.544"
$ws.Range("E8").Value = "  -0.81%  "

$ws.Range("D9").Value = "2.632.48"
$ws.Range("E9").Value = "  -0.11%  "

$ws.Range("D10").Value = "'0.136"
$ws.Range("E10").Value = "  +8.34%  "

$ws.Range("E11").Value = "  -0.57%  "

$ws.Range("D12").Value = "'5.20"
$ws.Range("E12").Value = "  -0.79%  "

$ws.Range("D13").Value = "'0.348"
$ws.Range("E13").Value = "  -1.51%  "

$ws.Range("B14").Value = "ShibaInu"
$ws.Range("C14").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D14").Value = "'0.0000191"
$ws.Range("E14").Value = "  +3.32%  "

$ws.Range("B15").Value = "Avalanche"
$ws.Range("C15").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D15").Value = "'27.58"
$ws.Range("E15").Value = "  -2.17%  "

$ws.Range("D16").Value = "3.120.00"
$ws.Range("E16").Value = "  +0.17%  "

$ws.Range("D17").Value = "67.873.49"
$ws.Range("E17").Value = "  +0.23%  "

$ws.Range("D18").Value = "2.634.92"
$ws.Range("E18").Value = "  -0.10%  "

$ws.Range("D19").Value = "'376.37"
$ws.Range("E19").Value = "  +2.27%  "

$ws.Range("D20").Value = "'11.18"
$ws.Range("E20").Value = "  -1.31%  "

$ws.Range("D21").Value = "'7.46"
$ws.Range("E21").Value = "  -2.24%  "

$ws.Range("D22").Value = "'4.23"
$ws.Range("E22").Value = "  -1.82%  "

$ws.Range("D23").Value = "'4.82"
$ws.Range("E23").Value = "  -3.41%  "

$ws.Range("D24").Value = "'2.04"
$ws.Range("E24").Value = "  -4.40%  "

$ws.Range("B25").Value = "Litecoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D25").Value = "'74.33"
$ws.Range("E25").Value = "  +5.19%  "

$ws.Range("B26").Value = "Dai"
$ws.Range("C26").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D26").Value = "'1.00"
$ws.Range("E26").Value = "  -0.01%  "

$ws.Range("B27").Value = "Aptos"
$ws.Range("C27").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D27").Value = "'9.88"
$ws.Range("E27").Value = "  -2.06%  "

$ws.Range("D28").Value = "2.776.94"
$ws.Range("E28").Value = "  +0.77%  "

$ws.Range("B29").Value = "PEPE"
$ws.Range("C29").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D29").Value = "'0.0000104"
$ws.Range("E29").Value = "  -0.52%  "

$ws.Range("D30").Value = "'589.27"
$ws.Range("E30").Value = "  -0.11%  "

$ws.Range("B31").Value = "Binance-PegBSC-USD"
$ws.Range("C31").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D31").Value = "'1.00"
$ws.Range("E31").Value = "  -0.07%  "

$ws.Range("B32").Value = "InternetComputer(DFINITY)"
$ws.Range("C32").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D32").Value = "'7.78"
$ws.Range("E32").Value = "  -2.07%  "

$ws.Range("B33").Value = "Fetch.AI"
$ws.Range("C33").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D33").Value = "'1.38"
$ws.Range("E33").Value = "  -3.98%  "

$ws.Range("D34").Value = "'1.84"
$ws.Range("E34").Value = "  -1.28%  "

$ws.Range("E35").Value = "  +0.01%  "

$ws.Range("D36").Value = "'0.127"
$ws.Range("E36").Value = "  -3.10%  "

$ws.Range("D37").Value = "'1.52"
$ws.Range("E37").Value = "  -1.22%  "

$ws.Range("D38").Value = "'158.57"
$ws.Range("E38").Value = "  +0.62%  "

$ws.Range("D39").Value = "'19.23"
$ws.Range("E39").Value = "  -1.77%  "

$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D40").Value = "'1.90"
$ws.Range("E40").Value = "  +2.96%  "

$ws.Range("B41").Value = "PolygonEcosystemToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D41").Value = "'0.368"
$ws.Range("E41").Value = "  -1.33%  "

$ws.Range("D42").Value = "'5.30"
$ws.Range("E42").Value = "  -1.28%  "

$ws.Range("D43").Value = "'2.64"
$ws.Range("E43").Value = "  -0.09%  "

$ws.Range("D44").Value = "'17.12"
$ws.Range("E44").Value = "  +4.47%  "

$ws.Range("D45").Value = "'1.00"
$ws.Range("E45").Value = "  +0.13%  "

$ws.Range("D46").Value = "'40.31"
$ws.Range("E46").Value = "  -2.26%  "

$ws.Range("D47").Value = "'155.75"
$ws.Range("E47").Value = "  -0.87%  "

$ws.Range("D48").Value = "0.0₆0293"
$ws.Range("E48").Value = "  +1.15%  "

$ws.Range("D49").Value = "'3.69"
$ws.Range("E49").Value = "  -1.69%  "

$ws.Range("D50").Value = "'1.70"
$ws.Range("E50").Value = "  -2.53%  "

$ws.Range("D51").Value = "'0.0781"
$ws.Range("E51").Value = "  -1.04%  "
